$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Islandora Metadata Template")

# --- Insert "Advisors" column before the current "Department" column (O) ---
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "Advisors"

# --- Insert "Embargo-Date" column before the current "Season" column (after Year) ---
# (Year is now at V1 after the first insert, Season moved to X1; new col goes at W)
$ws.Columns("W").Insert()
$ws.Range("W1").Value = "Embargo-Date"
# Give it its own text-formatted style (matches the "Year"/"Range" group fill,
# but stored as text rather than date) instead of inheriting the default.
$ws.Range("U1").Copy()
$ws.Range("W1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("W1").NumberFormat = "@"
$excel.CutCopyMode = $false

# --- Insert "Report-Number" column before the current "Rights" column ---
# (BePress-Identifier is now at BK1 after the previous inserts; new col goes at BL)
$ws.Columns("BL").Insert()
$ws.Range("BL1").Value = "Report-Number"

$wb.Save()
